$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ---- ALC ----
$wsALC.Cells.Item(64, 8).Value = 6311.96   # H64: 6287.1914 -> 6311.96
$wsALC.Cells.Item(64, 10).Value = 6425.5317   # J64: 6406.8184 -> 6425.5317
$wsALC.Cells.Item(64, 12).Value = 6425.5317   # L64: 6406.8184 -> 6425.5317
$wsALC.Cells.Item(64, 14).Value = -6921.5317   # N64: -6902.8184 -> -6921.5317
$wsALC.Cells.Item(67, 8).Value = 6311.96   # H67: 6287.1914 -> 6311.96
$wsALC.Cells.Item(67, 10).Value = 6425.5317   # J67: 6406.8184 -> 6425.5317
$wsALC.Cells.Item(67, 12).Value = 6425.5317   # L67: 6406.8184 -> 6425.5317
$wsALC.Cells.Item(67, 14).Value = -8141.5317   # N67: -8122.8184 -> -8141.5317
$wsALC.Cells.Item(86, 8).Value = 1951935.5   # H86: 2026991.1 -> 1951935.5
$wsALC.Cells.Item(86, 9).Value = 1513.2142   # I86: 1591.9231 -> 1513.2142
$wsALC.Cells.Item(86, 11).Value = 1513.2142   # K86: 1591.9231 -> 1513.2142
$wsALC.Cells.Item(86, 13).Value = -390.2141999999999   # M86: -468.9231 -> -390.2141999999999
$wsALC.Cells.Item(89, 8).Value = 1951935.5   # H89: 2026991.1 -> 1951935.5
$wsALC.Cells.Item(89, 9).Value = 1513.2142   # I89: 1591.9231 -> 1513.2142
$wsALC.Cells.Item(89, 11).Value = 7566.071   # K89: 7959.6155 -> 7566.071
$wsALC.Cells.Item(89, 13).Value = -1950.071   # M89: -2343.6155 -> -1950.071
$wsALC.Cells.Item(96, 8).Value = 626.5   # H96: 631.5 -> 626.5
$wsALC.Cells.Item(96, 9).Value = 605.9167   # I96: 607.8 -> 605.9167
$wsALC.Cells.Item(96, 11).Value = 1817.7501   # K96: 1823.4 -> 1817.7501
$wsALC.Cells.Item(96, 13).Value = -444.7501   # M96: -450.3999999999999 -> -444.7501
$wsALC.Cells.Item(111, 8).Value = 84615.164   # H111: 92223 -> 84615.164
$wsALC.Cells.Item(111, 9).Value = 112320.22   # I111: 126244.125 -> 112320.22
$wsALC.Cells.Item(111, 11).Value = 336960.66   # K111: 378732.375 -> 336960.66
$wsALC.Cells.Item(111, 13).Value = -333893.66   # M111: -375665.375 -> -333893.66
$wsALC.Cells.Item(132, 8).Value = 4425.933   # H132: 2364.3784 -> 4425.933
$wsALC.Cells.Item(132, 9).Value = 4782.75   # I132: 2337.7354 -> 4782.75
$wsALC.Cells.Item(132, 10).Value = 2998.6667   # J132: 2666.3333 -> 2998.6667
$wsALC.Cells.Item(132, 11).Value = 14348.25   # K132: 7013.206200000001 -> 14348.25
$wsALC.Cells.Item(132, 12).Value = 8996.000100000001   # L132: 7998.999899999999 -> 8996.000100000001
$wsALC.Cells.Item(132, 13).Value = -11818.25   # M132: -4483.206200000001 -> -11818.25
$wsALC.Cells.Item(132, 14).Value = -14056.0001   # N132: -13058.9999 -> -14056.0001
$wsALC.Cells.Item(135, 8).Value = 1671.9697   # H135: 1714.5312 -> 1671.9697
$wsALC.Cells.Item(135, 9).Value = 1621.7742   # I135: 1665.5 -> 1621.7742
$wsALC.Cells.Item(135, 11).Value = 14595.9678   # K135: 14989.5 -> 14595.9678
$wsALC.Cells.Item(135, 13).Value = -12060.9678   # M135: -12454.5 -> -12060.9678
$wsALC.Cells.Item(137, 8).Value = 1817.7693   # H137: 1839.7059 -> 1817.7693
$wsALC.Cells.Item(137, 9).Value = 1573.5143   # I137: 1598.0883 -> 1573.5143
$wsALC.Cells.Item(137, 10).Value = 2320.647   # J137: 2322.9412 -> 2320.647
$wsALC.Cells.Item(137, 11).Value = 4720.5429   # K137: 4794.2649 -> 4720.5429
$wsALC.Cells.Item(137, 12).Value = 6961.941   # L137: 6968.823600000001 -> 6961.941
$wsALC.Cells.Item(137, 13).Value = -2170.5429   # M137: -2244.2649 -> -2170.5429
$wsALC.Cells.Item(137, 14).Value = -12061.941   # N137: -12068.8236 -> -12061.941

# ---- ARM ----
$wsARM.Cells.Item(61, 8).Value = 4269.625   # H61: 4649.643 -> 4269.625
$wsARM.Cells.Item(61, 9).Value = 3954.2666   # I61: 4315 -> 3954.2666
$wsARM.Cells.Item(61, 11).Value = 3954.2666   # K61: 4315 -> 3954.2666
$wsARM.Cells.Item(61, 13).Value = -3742.2666   # M61: -4103 -> -3742.2666
$wsARM.Cells.Item(63, 8).Value = 5168.9165   # H63: 6752.0586 -> 5168.9165
$wsARM.Cells.Item(63, 9).Value = 2849   # I63: 4346.75 -> 2849
$wsARM.Cells.Item(63, 10).Value = 7910.636   # J63: 8890.111000000001 -> 7910.636
$wsARM.Cells.Item(63, 11).Value = 2849   # K63: 4346.75 -> 2849
$wsARM.Cells.Item(63, 12).Value = 7910.636   # L63: 8890.111000000001 -> 7910.636
$wsARM.Cells.Item(63, 13).Value = -2163   # M63: -3660.75 -> -2163
$wsARM.Cells.Item(63, 14).Value = -9282.636   # N63: -10262.111 -> -9282.636
$wsARM.Cells.Item(66, 8).Value = 5168.9165   # H66: 6752.0586 -> 5168.9165
$wsARM.Cells.Item(66, 9).Value = 2849   # I66: 4346.75 -> 2849
$wsARM.Cells.Item(66, 10).Value = 7910.636   # J66: 8890.111000000001 -> 7910.636
$wsARM.Cells.Item(66, 11).Value = 14245   # K66: 21733.75 -> 14245
$wsARM.Cells.Item(66, 12).Value = 39553.18   # L66: 44450.55500000001 -> 39553.18
$wsARM.Cells.Item(66, 13).Value = -10813   # M66: -18301.75 -> -10813
$wsARM.Cells.Item(66, 14).Value = -46417.18   # N66: -51314.55500000001 -> -46417.18
$wsARM.Cells.Item(82, 8).Value = 10385.333   # H82: 10718.667 -> 10385.333
$wsARM.Cells.Item(82, 10).Value = 10385.333   # J82: 10718.667 -> 10385.333
$wsARM.Cells.Item(82, 12).Value = 10385.333   # L82: 10718.667 -> 10385.333
$wsARM.Cells.Item(82, 14).Value = -11107.333   # N82: -11440.667 -> -11107.333
$wsARM.Cells.Item(85, 8).Value = 10385.333   # H85: 10718.667 -> 10385.333
$wsARM.Cells.Item(85, 10).Value = 10385.333   # J85: 10718.667 -> 10385.333
$wsARM.Cells.Item(85, 12).Value = 10385.333   # L85: 10718.667 -> 10385.333
$wsARM.Cells.Item(85, 14).Value = -12881.333   # N85: -13214.667 -> -12881.333
$wsARM.Cells.Item(136, 8).Value = 4269.625   # H136: 4649.643 -> 4269.625
$wsARM.Cells.Item(136, 9).Value = 3954.2666   # I136: 4315 -> 3954.2666
$wsARM.Cells.Item(136, 11).Value = 11862.7998   # K136: 12945 -> 11862.7998
$wsARM.Cells.Item(136, 13).Value = -9312.799800000001   # M136: -10395 -> -9312.799800000001
$wsARM.Cells.Item(139, 8).Value = 48101.57   # H139: 48928.75 -> 48101.57
$wsARM.Cells.Item(139, 10).Value = 48101.57   # J139: 48928.75 -> 48101.57
$wsARM.Cells.Item(139, 12).Value = 48101.57   # L139: 48928.75 -> 48101.57
$wsARM.Cells.Item(139, 14).Value = -58381.57   # N139: -59208.75 -> -58381.57

# ---- BSM ----
$wsBSM.Cells.Item(94, 8).Value = 2390.625   # H94: 2871.0833 -> 2390.625
$wsBSM.Cells.Item(94, 9).Value = 2896.0908   # I94: 3328.7778 -> 2896.0908
$wsBSM.Cells.Item(94, 10).Value = 1278.6   # J94: 1498 -> 1278.6
$wsBSM.Cells.Item(94, 11).Value = 2896.0908   # K94: 3328.7778 -> 2896.0908
$wsBSM.Cells.Item(94, 12).Value = 1278.6   # L94: 1498 -> 1278.6
$wsBSM.Cells.Item(94, 13).Value = -2445.0908   # M94: -2877.7778 -> -2445.0908
$wsBSM.Cells.Item(94, 14).Value = -2180.6   # N94: -2400 -> -2180.6

# ---- CRP ----
$wsCRP.Cells.Item(31, 8).Value = 64406.707   # H31: 76923.21000000001 -> 64406.707
$wsCRP.Cells.Item(31, 9).Value = 3058.5833   # I31: 3291.4 -> 3058.5833
$wsCRP.Cells.Item(31, 10).Value = 211642.2   # J31: 261002.75 -> 211642.2
$wsCRP.Cells.Item(31, 11).Value = 3058.5833   # K31: 3291.4 -> 3058.5833
$wsCRP.Cells.Item(31, 12).Value = 211642.2   # L31: 261002.75 -> 211642.2
$wsCRP.Cells.Item(31, 13).Value = -2763.5833   # M31: -2996.4 -> -2763.5833
$wsCRP.Cells.Item(31, 14).Value = -212232.2   # N31: -261592.75 -> -212232.2
$wsCRP.Cells.Item(34, 8).Value = 64406.707   # H34: 76923.21000000001 -> 64406.707
$wsCRP.Cells.Item(34, 9).Value = 3058.5833   # I34: 3291.4 -> 3058.5833
$wsCRP.Cells.Item(34, 10).Value = 211642.2   # J34: 261002.75 -> 211642.2
$wsCRP.Cells.Item(34, 11).Value = 3058.5833   # K34: 3291.4 -> 3058.5833
$wsCRP.Cells.Item(34, 12).Value = 211642.2   # L34: 261002.75 -> 211642.2
$wsCRP.Cells.Item(34, 13).Value = -2856.5833   # M34: -3089.4 -> -2856.5833
$wsCRP.Cells.Item(34, 14).Value = -212046.2   # N34: -261406.75 -> -212046.2

# ---- CUL ----
$wsCUL.Cells.Item(107, 8).Value = 56322.74   # H107: 59396.277 -> 56322.74
$wsCUL.Cells.Item(107, 10).Value = 88526.164   # J107: 96483.17999999999 -> 88526.164
$wsCUL.Cells.Item(107, 12).Value = 265578.492   # L107: 289449.54 -> 265578.492
$wsCUL.Cells.Item(107, 14).Value = -269418.492   # N107: -293289.54 -> -269418.492
$wsCUL.Cells.Item(132, 8).Value = 443093.72   # H132: 461480.97 -> 443093.72
$wsCUL.Cells.Item(132, 10).Value = 591944.2   # J132: 628828.2 -> 591944.2
$wsCUL.Cells.Item(132, 12).Value = 5327497.8   # L132: 5659453.8 -> 5327497.8
$wsCUL.Cells.Item(132, 14).Value = -5332557.8   # N132: -5664513.8 -> -5332557.8

# ---- GSM ----
$wsGSM.Cells.Item(51, 8).Value = 0   # H51: 14750 -> 0
$wsGSM.Cells.Item(51, 10).Value = 0   # J51: 14750 -> 0
$wsGSM.Cells.Item(51, 12).Value = 0   # L51: 14750 -> 0
$wsGSM.Cells.Item(51, 14).ClearContents()   # N51: delete (was -15768)
$wsGSM.Cells.Item(70, 8).Value = 11782.333   # H70: 12071.6 -> 11782.333
$wsGSM.Cells.Item(70, 9).Value = 8993.182000000001   # I70: 9292.799999999999 -> 8993.182000000001
$wsGSM.Cells.Item(70, 11).Value = 8993.182000000001   # K70: 9292.799999999999 -> 8993.182000000001
$wsGSM.Cells.Item(70, 13).Value = -8723.182000000001   # M70: -9022.799999999999 -> -8723.182000000001
$wsGSM.Cells.Item(73, 8).Value = 11782.333   # H73: 12071.6 -> 11782.333
$wsGSM.Cells.Item(73, 9).Value = 8993.182000000001   # I73: 9292.799999999999 -> 8993.182000000001
$wsGSM.Cells.Item(73, 11).Value = 8993.182000000001   # K73: 9292.799999999999 -> 8993.182000000001
$wsGSM.Cells.Item(73, 13).Value = -8057.182000000001   # M73: -8356.799999999999 -> -8057.182000000001
$wsGSM.Cells.Item(97, 8).Value = 1245.45   # H97: 1258.4736 -> 1245.45
$wsGSM.Cells.Item(97, 10).Value = 999.25   # J97: 999.6667 -> 999.25
$wsGSM.Cells.Item(97, 12).Value = 999.25   # L97: 999.6667 -> 999.25
$wsGSM.Cells.Item(97, 14).Value = -1991.25   # N97: -1991.6667 -> -1991.25
$wsGSM.Cells.Item(136, 8).Value = 26222.562   # H136: 27193.932 -> 26222.562
$wsGSM.Cells.Item(136, 10).Value = 26222.562   # J136: 27193.932 -> 26222.562
$wsGSM.Cells.Item(136, 12).Value = 78667.686   # L136: 81581.796 -> 78667.686
$wsGSM.Cells.Item(136, 14).Value = -83767.686   # N136: -86681.796 -> -83767.686
$wsGSM.Cells.Item(141, 8).Value = 132000   # H141: 133500 -> 132000
$wsGSM.Cells.Item(141, 10).Value = 132000   # J141: 133500 -> 132000
$wsGSM.Cells.Item(141, 12).Value = 132000   # L141: 133500 -> 132000
$wsGSM.Cells.Item(141, 14).Value = -142360   # N141: -143860 -> -142360

# ---- LTW ----
$wsLTW.Cells.Item(22, 8).Value = 499.5   # H22: 499 -> 499.5
$wsLTW.Cells.Item(22, 10).Value = 500   # J22: 0 -> 500
$wsLTW.Cells.Item(22, 12).Value = 500   # L22: 0 -> 500
$wsLTW.Cells.Item(22, 14).Value = -1090   # N22: None -> -1090
$wsLTW.Cells.Item(27, 8).Value = 499.5   # H27: 499 -> 499.5
$wsLTW.Cells.Item(27, 10).Value = 500   # J27: 0 -> 500
$wsLTW.Cells.Item(27, 12).Value = 500   # L27: 0 -> 500
$wsLTW.Cells.Item(27, 14).Value = -714   # N27: None -> -714
$wsLTW.Cells.Item(61, 8).Value = 9999.75   # H61: 12633 -> 9999.75
$wsLTW.Cells.Item(61, 9).Value = 16500   # I61: 30000 -> 16500
$wsLTW.Cells.Item(61, 10).Value = 3499.5   # J61: 3949.5 -> 3499.5
$wsLTW.Cells.Item(61, 11).Value = 16500   # K61: 30000 -> 16500
$wsLTW.Cells.Item(61, 12).Value = 3499.5   # L61: 3949.5 -> 3499.5
$wsLTW.Cells.Item(61, 13).Value = -16298   # M61: -29798 -> -16298
$wsLTW.Cells.Item(61, 14).Value = -3903.5   # N61: -4353.5 -> -3903.5
$wsLTW.Cells.Item(82, 8).Value = 1565.1177   # H82: 1553.2941 -> 1565.1177
$wsLTW.Cells.Item(82, 9).Value = 1403.8889   # I82: 1292.9 -> 1403.8889
$wsLTW.Cells.Item(82, 10).Value = 1746.5   # J82: 1925.2858 -> 1746.5
$wsLTW.Cells.Item(82, 11).Value = 1403.8889   # K82: 1292.9 -> 1403.8889
$wsLTW.Cells.Item(82, 12).Value = 1746.5   # L82: 1925.2858 -> 1746.5
$wsLTW.Cells.Item(82, 13).Value = -1042.8889   # M82: -931.9000000000001 -> -1042.8889
$wsLTW.Cells.Item(82, 14).Value = -2468.5   # N82: -2647.2858 -> -2468.5
$wsLTW.Cells.Item(85, 8).Value = 1565.1177   # H85: 1553.2941 -> 1565.1177
$wsLTW.Cells.Item(85, 9).Value = 1403.8889   # I85: 1292.9 -> 1403.8889
$wsLTW.Cells.Item(85, 10).Value = 1746.5   # J85: 1925.2858 -> 1746.5
$wsLTW.Cells.Item(85, 11).Value = 1403.8889   # K85: 1292.9 -> 1403.8889
$wsLTW.Cells.Item(85, 12).Value = 1746.5   # L85: 1925.2858 -> 1746.5
$wsLTW.Cells.Item(85, 13).Value = -155.8888999999999   # M85: -44.90000000000009 -> -155.8888999999999
$wsLTW.Cells.Item(85, 14).Value = -4242.5   # N85: -4421.2858 -> -4242.5
$wsLTW.Cells.Item(113, 8).Value = 9999.75   # H113: 12633 -> 9999.75
$wsLTW.Cells.Item(113, 9).Value = 16500   # I113: 30000 -> 16500
$wsLTW.Cells.Item(113, 10).Value = 3499.5   # J113: 3949.5 -> 3499.5
$wsLTW.Cells.Item(113, 11).Value = 16500   # K113: 30000 -> 16500
$wsLTW.Cells.Item(113, 12).Value = 3499.5   # L113: 3949.5 -> 3499.5
$wsLTW.Cells.Item(113, 13).Value = -14330   # M113: -27830 -> -14330
$wsLTW.Cells.Item(113, 14).Value = -7839.5   # N113: -8289.5 -> -7839.5
$wsLTW.Cells.Item(132, 8).Value = 4056   # H132: 3702.4827 -> 4056
$wsLTW.Cells.Item(132, 9).Value = 3867.2   # I132: 3606.3704 -> 3867.2
$wsLTW.Cells.Item(132, 11).Value = 11601.6   # K132: 10819.1112 -> 11601.6
$wsLTW.Cells.Item(132, 13).Value = -9071.599999999999   # M132: -8289.111199999999 -> -9071.599999999999

# ---- WVR ----
$wsWVR.Cells.Item(119, 8).Value = 0   # H119: 60000 -> 0
$wsWVR.Cells.Item(119, 10).Value = 0   # J119: 60000 -> 0
$wsWVR.Cells.Item(119, 12).Value = 0   # L119: 60000 -> 0
$wsWVR.Cells.Item(119, 14).ClearContents()   # N119: delete (was -69676)
$wsWVR.Cells.Item(141, 8).Value = 48518.332   # H141: 50000 -> 48518.332
$wsWVR.Cells.Item(141, 10).Value = 48518.332   # J141: 50000 -> 48518.332
$wsWVR.Cells.Item(141, 12).Value = 48518.332   # L141: 50000 -> 48518.332
$wsWVR.Cells.Item(141, 14).Value = -58878.332   # N141: -60360 -> -58878.332
